$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, shifting existing rows 130-229 down to 131-230.
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new data record.
$ws.Cells.Item(130, 1).Value = 7
$ws.Cells.Item(130, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(130, 3).Value = "Ñuble"
$ws.Cells.Item(130, 4).Value = 44957
$ws.Cells.Item(130, 5).Value = 16
$ws.Cells.Item(130, 6).Value = 100112045
$ws.Cells.Item(130, 7).Value = "Zapallo"
$ws.Cells.Item(130, 8).Value = "Camote"
$ws.Cells.Item(130, 9).Value = "1a (cosecha)"
$ws.Cells.Item(130, 10).Value = 400
$ws.Cells.Item(130, 11).Value = 450
$ws.Cells.Item(130, 12).Value = 500
$ws.Cells.Item(130, 13).Value = 475
$ws.Cells.Item(130, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(130, 15).Value = "Región del Maule"
$ws.Cells.Item(130, 16).Value = 475
$ws.Cells.Item(130, 17).Value = 1
$ws.Cells.Item(130, 18).Value = "Hortaliza"
